$d = $word.ActiveDocument

# Append a new paragraph "En tredje ændring" at the very end of the
# document body, after the existing last paragraph ("Lille ændring2")
# and before the section break.
$end = $d.Content
$end.Collapse(0)
$end.InsertAfter([char]13 + "En tredje ændring")
